$d = $word.ActiveDocument

$replacements = @(
    @("399÷5=79, 4", "261÷7=37, 2"),
    @("881÷2=440, 1", "469÷3=156, 1"),
    @("451÷7=64, 3", "794÷6=132, 2"),
    @("310÷9=34, 4", "514÷2=257, 0"),
    @("354÷4=88, 2", "771÷6=128, 3"),
    @("255÷2=127, 1", "514÷7=73, 3"),
    @("931÷6=155, 1", "971÷7=138, 5"),
    @("324÷7=46, 2", "196÷3=65, 1"),
    @("791÷7=113, 0", "357÷8=44, 5"),
    @("189÷6=31, 3", "591÷2=295, 1"),
    @("564÷6=94, 0", "731÷3=243, 2"),
    @("117÷8=14, 5", "283÷6=47, 1"),
    @("524÷6=87, 2", "404÷7=57, 5"),
    @("539÷5=107, 4", "106÷6=17, 4"),
    @("902÷4=225, 2", "180÷5=36, 0"),
    @("234÷6=39, 0", "104÷4=26, 0"),
    @("562÷4=140, 2", "503÷6=83, 5"),
    @("215÷4=53, 3", "368÷8=46, 0"),
    @("323÷6=53, 5", "275÷9=30, 5"),
    @("869÷5=173, 4", "476÷9=52, 8"),
    @("288÷3=96, 0", "536÷4=134, 0"),
    @("217÷4=54, 1", "720÷9=80, 0"),
    @("292÷9=32, 4", "840÷8=105, 0"),
    @("538÷9=59, 7", "792÷6=132, 0"),
    @("466÷6=77, 4", "578÷9=64, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
